# Actualización automática 2025-11-10 17:30:08
$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("D14").Value = 457.92
$ws1.Range("I14").Value = 44.1
$ws1.Range("L14").Value = 1165.07
$ws1.Range("M14").Value = 657.72

$ws1.Range("D23").Value = "1 de 21"
$ws1.Range("I23").Value = "1 de 21"
$ws1.Range("L23").Value = "1 de 21"
$ws1.Range("M23").Value = "3 de 21"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F14").Value = 2324.81
$ws2.Range("F23").Value = 4066.3

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D3").Value = 305.28
$ws3.Range("E3").Value = 2258.72
$ws3.Range("F3").Value = 0.1190639625585023

$ws3.Range("D7").Value = 44.1
$ws3.Range("E7").Value = 339.3
$ws3.Range("F7").Value = 0.1150234741784038

$ws3.Range("D11").Value = 13.67
$ws3.Range("E11").Value = 2664.33
$ws3.Range("F11").Value = 0.005104555638536221

$ws3.Range("D12").Value = 3703.25
$ws3.Range("E12").Value = 40714.75
$ws3.Range("F12").Value = 0.08337273177540637

$ws3.Range("D14").Value = 4066.3
$ws3.Range("E14").Value = 51333.17101170094
$ws3.Range("F14").Value = 0.07339961782561345
